# Applies the "improve checks for all Edm Types" edit to TwoRowsErrors.xlsx:
#  - fixes two Quantity values (B6, B9) to 2.5, B6 additionally formatted as Text
#  - replaces B10's numeric Quantity with the text value "2 EUR"
#  - corrects the computed "time" column for rows 15/16
#  - appends three new data rows (17-19) plus two blank formatted rows (20-21)
#  - widens columns C, G and I to fit the new content
#  - moves the active selection to H16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix existing rows -----------------------------------------------------
$ws.Range("B6").Value = 2.5
$ws.Range("B6").NumberFormat = "@"

$ws.Range("B9").Value = 2.5

$ws.Range("B10").Value = "2 EUR"

$ws.Range("H15").Value = 0.041666666666666664
$ws.Range("H16").Value = 0.95833333333333337

# --- append new data rows ----------------------------------------------------
# Pre-apply the date/time number formats used elsewhere in the sheet (columns
# E-G = short date, H = time) to the new rows by copying the formatting from
# row 2, then fill in the values/text for each cell individually.
$ws.Range("E2:G2").Copy()
$ws.Range("E17:G19").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H17:H18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A17").Value = 253
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = "Product Test 16"
$ws.Range("D17").Value = 27.56
$ws.Range("G17").Value = "23.11.2023456"
$ws.Range("F17").Value = 45254
$ws.Range("H17").Value = 0

$ws.Range("A18").Value = 253
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "Product Test 17"
$ws.Range("D18").Value = 28.56
$ws.Range("E18").Value = 45255
$ws.Range("G18").Value = "23.11.2023457"
$ws.Range("H18").Value = 1.2916666666666701

$ws.Range("A19").Value = 253
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "Product Test 18"
$ws.Range("D19").Value = 29.56
$ws.Range("E19").Value = 45255
$ws.Range("F19").Value = 45254
$ws.Range("G19").Value = "23.11.2023458"
$ws.Range("H19").Value = 2.041666666666667

# these two text values were corrected last (matches the original authoring
# order, which is also why they land at the end of the shared-string table)
$ws.Range("E17").Value = "25.11.202312"
$ws.Range("F18").Value = "24.11.202312"

# two further blank rows, formatted the same way as the data rows above
$ws.Range("E2:G2").Copy()
$ws.Range("E20:G21").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H20:H21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- column widths -----------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 12.6666666666667
$ws.Columns.Item(7).ColumnWidth = 15.1666666666667
$ws.Columns.Item(9).ColumnWidth = 17.0833333333333

# --- selection ---------------------------------------------------------------
$ws.Range("H16").Select()
